$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 462.42856
$ws.Range("I33").Value = 462.42856
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 462.42856
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -233.42856
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 77693.75
$ws.Range("I40").Value = 78591.664
$ws.Range("J40").Value = 75000
$ws.Range("K40").Value = 78591.664
$ws.Range("L40").Value = 75000
$ws.Range("M40").Value = -78416.664
$ws.Range("N40").Value = -75350
$ws.Range("H112").Value = 4392.9653
$ws.Range("I112").Value = 1350
$ws.Range("J112").Value = 4618.3706
$ws.Range("K112").Value = 4050
$ws.Range("L112").Value = 13855.1118
$ws.Range("M112").Value = -2942
$ws.Range("N112").Value = -16071.1118
$ws.Range("H132").Value = 187262.5
$ws.Range("I132").Value = 234235.88
$ws.Range("K132").Value = 702707.64
$ws.Range("M132").Value = -700177.64
$ws.Range("H137").Value = 529110.6
$ws.Range("I137").Value = 590758.3
$ws.Range("K137").Value = 1772274.9
$ws.Range("M137").Value = -1769724.9
$ws.Range("H138").Value = 5942.5225
$ws.Range("I138").Value = 1699.45
$ws.Range("J138").Value = 7748.085
$ws.Range("K138").Value = 5098.35
$ws.Range("L138").Value = 23244.255
$ws.Range("M138").Value = 41.64999999999964
$ws.Range("N138").Value = -33524.255
$ws.Range("H141").Value = 2631.739
$ws.Range("I141").Value = 2415.1775
$ws.Range("J141").Value = 4549.857
$ws.Range("K141").Value = 7245.532499999999
$ws.Range("L141").Value = 13649.571
$ws.Range("M141").Value = -2065.532499999999
$ws.Range("N141").Value = -24009.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15048.378
$ws.Range("I32").Value = 15826.709
$ws.Range("J32").Value = 12795.315
$ws.Range("K32").Value = 15826.709
$ws.Range("L32").Value = 12795.315
$ws.Range("M32").Value = -15539.709
$ws.Range("N32").Value = -13369.315
$ws.Range("H61").Value = 3882.796
$ws.Range("I61").Value = 2597.7942
$ws.Range("J61").Value = 6795.467
$ws.Range("K61").Value = 2597.7942
$ws.Range("L61").Value = 6795.467
$ws.Range("M61").Value = -2385.7942
$ws.Range("N61").Value = -7219.467
$ws.Range("H132").Value = 18488.764
$ws.Range("I132").Value = 23976.27
$ws.Range("J132").Value = 6599.1665
$ws.Range("K132").Value = 71928.81
$ws.Range("L132").Value = 19797.4995
$ws.Range("M132").Value = -69398.81
$ws.Range("N132").Value = -24857.4995
$ws.Range("H136").Value = 3882.796
$ws.Range("I136").Value = 2597.7942
$ws.Range("J136").Value = 6795.467
$ws.Range("K136").Value = 7793.382599999999
$ws.Range("L136").Value = 20386.401
$ws.Range("M136").Value = -5243.382599999999
$ws.Range("N136").Value = -25486.401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3946
$ws.Range("I105").Value = 4050.6
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 4050.6
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -2303.6
$ws.Range("N105").Value = -6394
$ws.Range("H134").Value = 2037.2778
$ws.Range("I134").Value = 1746.6207
$ws.Range("J134").Value = 3241.4285
$ws.Range("K134").Value = 5239.8621
$ws.Range("L134").Value = 9724.2855
$ws.Range("M134").Value = -2704.8621
$ws.Range("N134").Value = -14794.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H31").Value = 2749.1345
$ws.Range("I31").Value = 1465.186
$ws.Range("J31").Value = 8883.556
$ws.Range("K31").Value = 1465.186
$ws.Range("L31").Value = 8883.556
$ws.Range("M31").Value = -1170.186
$ws.Range("N31").Value = -9473.556
$ws.Range("H34").Value = 2749.1345
$ws.Range("I34").Value = 1465.186
$ws.Range("J34").Value = 8883.556
$ws.Range("K34").Value = 1465.186
$ws.Range("L34").Value = 8883.556
$ws.Range("M34").Value = -1263.186
$ws.Range("N34").Value = -9287.556
$ws.Range("H58").Value = 2000.3334
$ws.Range("I58").Value = 1794.1482
$ws.Range("J58").Value = 2618.889
$ws.Range("K58").Value = 1794.1482
$ws.Range("L58").Value = 2618.889
$ws.Range("M58").Value = -1591.1482
$ws.Range("N58").Value = -3024.889
$ws.Range("H62").Value = 49500.855
$ws.Range("I62").Value = 6750
$ws.Range("J62").Value = 66601.2
$ws.Range("K62").Value = 6750
$ws.Range("L62").Value = 66601.2
$ws.Range("M62").Value = -6126
$ws.Range("N62").Value = -67849.2
$ws.Range("H65").Value = 49500.855
$ws.Range("I65").Value = 6750
$ws.Range("J65").Value = 66601.2
$ws.Range("K65").Value = 33750
$ws.Range("L65").Value = 333006
$ws.Range("M65").Value = -30630
$ws.Range("N65").Value = -339246
$ws.Range("H99").Value = 25038.309
$ws.Range("I99").Value = 49874.5
$ws.Range("J99").Value = 14000
$ws.Range("K99").Value = 49874.5
$ws.Range("L99").Value = 14000
$ws.Range("M99").Value = -48376.5
$ws.Range("N99").Value = -16996
$ws.Range("H105").Value = 1200.4546
$ws.Range("I105").Value = 815.6667
$ws.Range("J105").Value = 1662.2
$ws.Range("K105").Value = 815.6667
$ws.Range("L105").Value = 1662.2
$ws.Range("M105").Value = 931.3333
$ws.Range("N105").Value = -5156.2
$ws.Range("H107").Value = 808.2857
$ws.Range("I107").Value = 677.3077
$ws.Range("J107").Value = 1021.125
$ws.Range("K107").Value = 677.3077
$ws.Range("L107").Value = 1021.125
$ws.Range("M107").Value = 1242.6923
$ws.Range("N107").Value = -4861.125
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 25038.309
$ws.Range("I126").Value = 49874.5
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 149623.5
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -147153.5
$ws.Range("N126").Value = -46940
$ws.Range("H134").Value = 1996.2142
$ws.Range("I134").Value = 1779.7576
$ws.Range("J134").Value = 2789.889
$ws.Range("K134").Value = 5339.2728
$ws.Range("L134").Value = 8369.667000000001
$ws.Range("M134").Value = -2804.2728
$ws.Range("N134").Value = -13439.667
$ws.Range("H136").Value = 2000.3334
$ws.Range("I136").Value = 1794.1482
$ws.Range("J136").Value = 2618.889
$ws.Range("K136").Value = 5382.444600000001
$ws.Range("L136").Value = 7856.667
$ws.Range("M136").Value = -2832.444600000001
$ws.Range("N136").Value = -12956.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 765.4
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 765.4
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2296.2
$ws.Range("N107").Value = -6136.2
$ws.Range("M107").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 36139.727
$ws.Range("J123").Value = 36139.727
$ws.Range("L123").Value = 36139.727
$ws.Range("N123").Value = -41039.727
$ws.Range("H140").Value = 58376.43
$ws.Range("J140").Value = 58376.43
$ws.Range("L140").Value = 58376.43
$ws.Range("N140").Value = -68736.42999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1099.5
$ws.Range("I100").Value = 1074.8334
$ws.Range("J100").Value = 1247.5
$ws.Range("K100").Value = 1074.8334
$ws.Range("L100").Value = 1247.5
$ws.Range("M100").Value = -533.8334
$ws.Range("N100").Value = -2329.5
$ws.Range("H136").Value = 3136.9062
$ws.Range("I136").Value = 2223.3408
$ws.Range("K136").Value = 6670.0224
$ws.Range("M136").Value = -4120.0224

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1120.875
$ws.Range("I132").Value = 799.0161000000001
$ws.Range("J132").Value = 2229.5
$ws.Range("K132").Value = 2397.0483
$ws.Range("L132").Value = 6688.5
$ws.Range("M132").Value = 132.9516999999996
$ws.Range("N132").Value = -11748.5
$ws.Range("H136").Value = 5381.2
$ws.Range("I136").Value = 1025.9215
$ws.Range("J136").Value = 9914.245000000001
$ws.Range("K136").Value = 3077.7645
$ws.Range("L136").Value = 29742.735
$ws.Range("M136").Value = -527.7644999999998
$ws.Range("N136").Value = -34842.735
